$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the plate number on row 2 (A2) first so the freed shared-string
# slot gets reused, matching how the edits were actually entered.
$ws.Range("A2").Value = "SZV190"

# Duplicate row 2 into row 3 so formatting/styles/number types carry over,
# then overwrite the cells that actually differ for the new record.
$ws.Range("A2:S2").Copy($ws.Range("A3:S3"))

# New contact record values for row 3.
$ws.Range("R3").Value = "tes2t@gmail.com"
$ws.Range("M3").Value = "Calle 1234"
$ws.Range("A3").Value = "GBL76F"
$ws.Range("O3").Value = 6664331
$ws.Range("S3").Value = 3054665669

# Extend the existing list-based data validations to cover the new row.
$ws.Range("N2").Validation.Delete()
$ws.Range("B2").Validation.Delete()
$ws.Range("N2:N3").Validation.Add(3, 1, 1, '"NA, BEEPER, CASA, CELULAR, E-MAIL, OFICINA, OFICINA 2, OFICINA 3, PBX / CONMUTADOR, RESIDENCIA, RESIDENCIA 2, TELEFAX, TELFAX 2, TRANSFERENCIAS"')
$ws.Range("B2:B3").Validation.Add(3, 1, 1, '"NA,CEDULA DE CIUDADANIA, CEDULA DE EXTRANJERIA, NRO DE NIT, PASAPORTE, TARJETA DE IDENTIDAD"')

# Add the mailto hyperlink for the new e-mail address (app URL update).
[void]$ws.Hyperlinks.Add($ws.Range("R3"), "mailto:tes2t@gmail.com")
$ws.Range("R3").Style = $ws.Range("R2").Style

# Update the selection shown in the sheet view.
[void]$ws.Range("B2:C7").Select()
